$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204659104347229
$ws.Range("B1").Value = 2.286876678466797
$ws.Range("C1").Value = 6.253408432006836
$ws.Range("D1").Value = 2.01954460144043
$ws.Range("E1").Value = 1.17466139793396
